$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a1"
$ws.Range("C2").Value = "Ddr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 3.048502333333333
$ws.Range("N2").Value = 9.145507
$ws.Range("O2").Value = 0.02767295150267144
$ws.Range("P2").Value = 0.02767295150267144
$ws.Range("Q2").Value = 18.51090552180566
$ws.Range("R2").Value = 166.598149696251
$ws.Range("S2").Value = 0.0001091421004432715
$ws.Range("T2").Value = 0.0001091421004432715

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a1"
$ws.Range("C3").Value = "Ddr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 83.828771
$ws.Range("N3").Value = 251.486313
$ws.Range("O3").Value = 0.7609603866942148
$ws.Range("P3").Value = 0.7609603866942148
$ws.Range("Q3").Value = 509.019279081001
$ws.Range("R3").Value = 4581.173511729009
$ws.Range("S3").Value = 0.00300122720736576
$ws.Range("T3").Value = 0.00300122720736576

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a1"
$ws.Range("C4").Value = "Ddr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.072131
$ws.Range("H4").Value = 18.216393
$ws.Range("I4").Value = 0.003943999267036455
$ws.Range("J4").Value = 0.003943999267036454
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 23.28453333333333
$ws.Range("N4").Value = 69.8536
$ws.Range("O4").Value = 0.2113666618031137
$ws.Range("P4").Value = 0.2113666618031138
$ws.Range("Q4").Value = 141.3867366738667
$ws.Range("R4").Value = 1272.4806300648
$ws.Range("S4").Value = 0.0008336299592274228
$ws.Range("T4").Value = 0.0008336299592274227

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a1"
$ws.Range("C5").Value = "Ddr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 3.048502333333333
$ws.Range("N5").Value = 9.145507
$ws.Range("O5").Value = 0.02767295150267144
$ws.Range("P5").Value = 0.02767295150267144
$ws.Range("Q5").Value = 4514.380187944215
$ws.Range("R5").Value = 40629.42169149794
$ws.Range("S5").Value = 0.02661722492891102
$ws.Range("T5").Value = 0.02661722492891101

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a1"
$ws.Range("C6").Value = "Ddr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1480.851806666667
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9618498744646554
$ws.Range("J6").Value = 0.9618498744646552
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 83.828771
$ws.Range("N6").Value = 251.486313
$ws.Range("O6").Value = 0.7609603866942148
$ws.Range("P6").Value = 0.7609603866942148
$ws.Range("Q6").Value = 124137.9869859963
$ws.Range("R6").Value = 1117241.882873966
$ws.Range("S6").Value = 0.7319296524144062
$ws.Range("T6").Value = 0.731929652414406

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a1"
$ws.Range("C7").Value = "Ddr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1480.851806666667
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9618498744646554
$ws.Range("J7").Value = 0.9618498744646552
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 23.28453333333333
$ws.Range("N7").Value = 69.8536
$ws.Range("O7").Value = 0.2113666618031137
$ws.Range("P7").Value = 0.2113666618031138
$ws.Range("Q7").Value = 34480.94325405688
$ws.Range("R7").Value = 310328.489286512
$ws.Range("S7").Value = 0.2033029971213382
$ws.Range("T7").Value = 0.2033029971213382

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col1a1"
$ws.Range("C8").Value = "Ddr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 52.663316
$ws.Range("H8").Value = 157.989948
$ws.Range("I8").Value = 0.03420612626830831
$ws.Range("J8").Value = 0.0342061262683083
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 3.048502333333333
$ws.Range("N8").Value = 9.145507
$ws.Range("O8").Value = 0.02767295150267144
$ws.Range("P8").Value = 0.02767295150267144
$ws.Range("Q8").Value = 160.5442417070707
$ws.Range("R8").Value = 1444.898175363636
$ws.Range("S8").Value = 0.0009465844733171513
$ws.Range("T8").Value = 0.0009465844733171512

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col1a1"
$ws.Range("C9").Value = "Ddr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 52.663316
$ws.Range("H9").Value = 157.989948
$ws.Range("I9").Value = 0.03420612626830831
$ws.Range("J9").Value = 0.0342061262683083
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 83.828771
$ws.Range("N9").Value = 251.486313
$ws.Range("O9").Value = 0.7609603866942148
$ws.Range("P9").Value = 0.7609603866942148
$ws.Range("Q9").Value = 4414.701057064636
$ws.Range("R9").Value = 39732.30951358172
$ws.Range("S9").Value = 0.02602950707244303
$ws.Range("T9").Value = 0.02602950707244303

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col1a1"
$ws.Range("C10").Value = "Ddr2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 52.663316
$ws.Range("H10").Value = 157.989948
$ws.Range("I10").Value = 0.03420612626830831
$ws.Range("J10").Value = 0.0342061262683083
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.28453333333333
$ws.Range("N10").Value = 69.8536
$ws.Range("O10").Value = 0.2113666618031137
$ws.Range("P10").Value = 0.2113666618031138
$ws.Range("Q10").Value = 1226.240736845867
$ws.Range("R10").Value = 11036.1666316128
$ws.Range("S10").Value = 0.007230034722548127
$ws.Range("T10").Value = 0.007230034722548126
